# Convert the title "Untitled" -> "Coursera Assignment" and the
# author "Muskaan" -> "Muskaan Dudeja", each split across three
# separate runs (word / space / word), matching an R-Markdown-style
# re-render of the heading + author line.

$d = $word.ActiveDocument

$p2 = $d.Paragraphs(2)

# Range spanning from the very start of the document through the end
# of the author paragraph (excluding its paragraph mark).
$rAll = $d.Range(0, $p2.Range.End - 1)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Coursera</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Assignment</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Muskaan</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Dudeja</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$rAll.InsertXML($xml)
